$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scoreboard rows for week of 6/25/2024 (date serial 45468)
$newRows = @(
    @{ Row=108; A="Jeremiah"; C="Workout";  D=26; E=0;    F=0;   G=24; H=2;  I=0;  J=0; K=0; L="Agile Antelope"; M=3 },
    @{ Row=109; A="Matt";     C="Run";      D=48; E=4.69; F=115; G=3;  H=25; I=11; J=3; K=4; L="Agile Antelope"; M=3 },
    @{ Row=110; A="Steven";   C="Run";      D=36; E=3.12; F=95;  G=1;  H=14; I=17; J=3; K=0; L="Wily Hyena";     M=3 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = 45468
    $ws.Cells.Item($row, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Match the post-edit selection/active cell recorded in the workbook
$ws.Range("A111").Select() | Out-Null
